# Change test case type (TC_KIND) from SCRIPTED to GHERKIN, which
# effectively removes the now-unused "SCRIPTED" shared string and adds
# the language info (TC_SCRIPTING_LANGUAGE already held GHERKIN).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST_CASES")
$ws.Activate()

# X2 holds the TC_KIND value, currently "SCRIPTED" -> change to "GHERKIN"
$ws.Range("X2").Value = "GHERKIN"

# Update the view's active cell/selection to X3, matching the recorded
# workbook state after the edit.
$ws.Range("X3").Select()
